$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23, pushing existing rows 23-24 down to 24-25
$ws.Rows.Item(23).Insert()

# Fill in new row 23 with the inserted week's data
$ws.Cells.Item(23, 1).Value = 1
$ws.Cells.Item(23, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(23, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(23, 4).Value = 44943
$ws.Cells.Item(23, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23, 5).Value = 15
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100103
$ws.Cells.Item(23, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(23, 9).Value = 100103001
$ws.Cells.Item(23, 10).Value = "Cereza"
$ws.Cells.Item(23, 11).Value = "Santina"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 600
$ws.Cells.Item(23, 14).Value = 14000
$ws.Cells.Item(23, 15).Value = 15000
$ws.Cells.Item(23, 16).Value = 14333
$ws.Cells.Item(23, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(23, 18).Value = "Región del Maule"
$ws.Cells.Item(23, 19).Value = 956
$ws.Cells.Item(23, 20).Value = 15
